# Saudi Arabia Division 1 - base update (06-04-2024 15:39)
#
# The upstream scraper re-sorted a handful of fixtures that share the exact
# same kick-off timestamp, swapping their row order. Net effect: for each
# such pair of rows, the whole payload (every column except the running
# "id" in column A) moves from one row to the other. Everything else in
# the sheet (including all other "Al Orubah"/"Al Jabalain" fixtures that
# are untouched by the reorder) keeps its original value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC == numeric column indices 2..29. Column A (the running
# "id") is left untouched since it simply enumerates row order top-to-bottom.
$firstCol = 2
$lastCol = 29

# Row pairs whose entire payload (B..AC) got exchanged between the two
# rows (same match date/time, order flipped in the refreshed export).
$fullSwapPairs = @(
    @(2, 3),
    @(49, 50),
    @(59, 60),
    @(68, 69),
    @(115, 116),
    @(128, 129),
    @(172, 173),
    @(203, 204),
    @(212, 213),
    @(238, 239),
    @(241, 242)
)

foreach ($pair in $fullSwapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $vals1 = @()
    $vals2 = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals1 += , ($ws.Cells.Item($r1, $c).Value2)
        $vals2 += , ($ws.Cells.Item($r2, $c).Value2)
    }
    for ($i = 0; $i -lt $vals1.Count; $i++) {
        $c = $firstCol + $i
        $ws.Cells.Item($r1, $c).Value = $vals2[$i]
        $ws.Cells.Item($r2, $c).Value = $vals1[$i]
    }
}
